$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns so numeric-looking
# strings (e.g. "597.78") are not auto-converted to numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '68.239.88'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '3.711.10'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '597.78'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').Value = '167.61'
$ws.Range('E6').Value = '  -3.06%  '
$ws.Range('D7').Value = '3.706.78'
$ws.Range('E7').Value = '  -3.02%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').Value = '0.168'
$ws.Range('E10').Value = '  +5.46%  '
$ws.Range('D11').Value = '6.24'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  -2.01%  '
$ws.Range('D13').Value = '38.10'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').Value = '4.329.28'
$ws.Range('E15').Value = '  -3.06%  '
$ws.Range('D16').Value = '3.710.26'
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('D17').Value = '68.233.61'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '7.29'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('D20').Value = '17.18'
$ws.Range('E20').Value = '  +7.00%  '
$ws.Range('D21').Value = '490.06'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '9.22'
$ws.Range('E22').Value = '  -1.22%  '
$ws.Range('D23').Value = '0.721'
$ws.Range('E23').Value = '  -2.42%  '
$ws.Range('D24').Value = '84.73'
$ws.Range('E24').Value = '  -1.58%  '
$ws.Range('D25').Value = '0.0000142'
$ws.Range('E25').Value = '  +2.55%  '
$ws.Range('D26').Value = '2.30'
$ws.Range('E26').Value = '  -3.30%  '
$ws.Range('D27').Value = '12.25'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '10.08'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '2.92'
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('D31').Value = '7.87'
$ws.Range('E31').Value = '  +3.06%  '
$ws.Range('E32').Value = '  -2.69%  '
$ws.Range('D33').Value = '31.44'
$ws.Range('E33').Value = '  -5.16%  '
$ws.Range('D34').Value = '3.853.66'
$ws.Range('E34').Value = '  -2.91%  '
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('D36').Value = '3.655.71'
$ws.Range('E36').Value = '  -3.01%  '
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('E38').Value = '  -0.91%  '
$ws.Range('D39').Value = '5.81'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -3.83%  '
$ws.Range('D41').Value = '0.322'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').Value = '430.78'
$ws.Range('E42').Value = '  -4.68%  '
$ws.Range('D43').Value = '48.68'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('D44').Value = '1.96'
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('D45').Value = '2.84'
$ws.Range('E45').Value = '  -2.27%  '
$ws.Range('D46').Value = '8.41'
$ws.Range('E46').Value = '  +1.02%  '
$ws.Range('D48').Value = '40.30'
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('D49').Value = '141.02'
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('D50').Value = '2.760.21'
$ws.Range('E50').Value = '  -3.39%  '
$ws.Range('D51').Value = '0.0350'
$ws.Range('E51').Value = '  -0.99%  '
